$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("O2").Value = 0.7806153008439968
$ws.Range("P2").Value = 0.7806153008439968
$ws.Range("S2").Value = 0.7806153008439968
$ws.Range("T2").Value = 0.7806153008439968

# Row 3 updates
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.01343
$ws.Range("N3").Value = 0.04029
$ws.Range("O3").Value = 0.2193846991560033
$ws.Range("P3").Value = 0.2193846991560033
$ws.Range("Q3").Value = 0.01418497640333333
$ws.Range("R3").Value = 0.12766478763
$ws.Range("S3").Value = 0.2193846991560033
$ws.Range("T3").Value = 0.2193846991560033
